$wb = $excel.ActiveWorkbook

$metaWs = $wb.Worksheets.Item("Metadata")
$includeWs = $wb.Worksheets.Item("Include from Medication item ")

# Rename the "Include" worksheet
$includeWs.Name = "Include #0"

# Update the compiled "Date" metadata value
$metaWs.Range("B8").Value = "2024-09-12T14:01:50+00:00"

# Insert a new "Jurisdiction" metadata row above "Description" (row 11),
# matching the formatting of the existing metadata rows.
$metaWs.Rows.Item(11).Insert()
$metaWs.Range("A10:B10").Copy()
$metaWs.Range("A11:B11").PasteSpecial(-4122)
$metaWs.Cells.Item(11, 1).Value = "Jurisdiction"
$metaWs.Cells.Item(11, 2).Value = ""
